$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New negative_feedback rows (349-416) appended to the dataset
$data = @(
    @('negative_feedback','Terrible job!','I''m sorry to hear that.'),
    @('negative_feedback','You are awful!','I apologize for any inconvenience caused.'),
    @('negative_feedback','Terrible work!','I understand your disappointment.'),
    @('negative_feedback','You''re a failure!','I''ll do my best to improve.'),
    @('negative_feedback','Poorly done!','Mistakes happen, and I''ll learn from them.'),
    @('negative_feedback','Disastrous!','I''ll make sure to work on it.'),
    @('negative_feedback','You''re bombing it!','I''m here to help, so please let me know how I can assist you.'),
    @('negative_feedback','Terrible job!','I''ll take your feedback into consideration.'),
    @('negative_feedback','You missed it!','I''ll strive to do better next time.'),
    @('negative_feedback','You''re the worst!','I appreciate your feedback, and I''ll make the necessary adjustments.'),
    @('negative_feedback','Awful!','I understand your frustration.'),
    @('negative_feedback','You''re a dunce!','I''ll make sure to address this issue.'),
    @('negative_feedback','Horrible work!','I apologize for the inconvenience caused.'),
    @('negative_feedback','You''re unimpressive!','I''ll make sure to improve my performance.'),
    @('negative_feedback','Flawed!','I''ll work on my skills to provide better results.'),
    @('negative_feedback','Terrible job!','I''m sorry if I haven''t met your expectations.'),
    @('negative_feedback','You''re an amateur!','I''ll take note of your feedback and make improvements.'),
    @('negative_feedback','Mediocre work!','I''ll work on delivering better outcomes in the future.'),
    @('negative_feedback','You''re a nobody!','I apologize for any dissatisfaction caused.'),
    @('negative_feedback','Unimpressive!','I''ll take your comments into account and make necessary improvements.'),
    @('negative_feedback','You''re a disaster!','I''ll strive to exceed your expectations in the future.'),
    @('negative_feedback','Unimpressive work!','I''m sorry for any inconvenience. Your feedback is valuable to me.'),
    @('negative_feedback','You''re a loser!','I understand that this is not what you expected.'),
    @('negative_feedback','You made a mess!','I''ll make sure to enhance my performance in the future.'),
    @('negative_feedback','You''re a failure!','I apologize for any negative experience you''ve had.'),
    @('negative_feedback','You''re a fake!','I''ll do my best to rectify the situation.'),
    @('negative_feedback','Terrible work!','I''m sorry if my responses have been inadequate.'),
    @('negative_feedback','You''re bombing it today!','I''ll make sure to improve my understanding and provide more accurate answers.'),
    @('negative_feedback','You''re not gifted!','I''m sorry if I haven''t met your expectations. I''ll take your feedback on board.'),
    @('negative_feedback','You''re making no effort!','I understand that my response may not have been what you were looking for. I''ll strive to do better.'),
    @('negative_feedback','You''re a fraud!','I apologize for any disappointment caused. I''ll make sure to improve my performance.'),
    @('negative_feedback','You''re a failure!','I''m sorry if I haven''t been able to assist you effectively. I''ll make the necessary improvements.'),
    @('negative_feedback','You''re not a pro!','I apologize for any inconvenience or dissatisfaction caused. I''ll work on it.'),
    @('negative_feedback','You''re average!','I understand that my responses have fallen short. I''ll make the necessary adjustments.'),
    @('negative_feedback','You''re a disappointment!','I''m sorry if I haven''t been able to meet your needs. I''ll strive to improve my performance.'),
    @('negative_feedback','You''re not a master!','I apologize for any shortcomings. I''ll make sure to enhance my capabilities.'),
    @('negative_feedback','You did a terrible job!','I understand that my answers haven''t been satisfactory. I''ll make the necessary improvements.'),
    @('negative_feedback','You''re not a true professional!','I appreciate your feedback, and I''ll make the necessary adjustments.'),
    @('negative_feedback','Terrible work!','I understand your frustration.'),
    @('negative_feedback','You''re struggling!','I''ll make sure to address this issue.'),
    @('negative_feedback','You''re not talented!','I apologize for the inconvenience caused.'),
    @('negative_feedback','You''re not making any progress!','I''ll make sure to improve my performance.'),
    @('negative_feedback','You''re a fool!','I apologize if my answers haven''t been up to your standards. I''ll work on it.'),
    @('negative_feedback','You''re a liability!','I''m sorry if I haven''t met your expectations.'),
    @('negative_feedback','You''re an underachiever!','I''ll take note of your feedback and make improvements.'),
    @('negative_feedback','You''re not creative!','I''ll work on delivering better outcomes in the future.'),
    @('negative_feedback','You''re not productive!','I apologize for any dissatisfaction caused.'),
    @('negative_feedback','You''re not a wizard!','I''ll take your comments into account and make necessary improvements.'),
    @('negative_feedback','You''re not worst','I''ll strive to exceed your expectations in the future.'),
    @('negative_feedback','You''re not a genius!','I''m sorry for any inconvenience. Your feedback is valuable to me.'),
    @('negative_feedback','You''re a weak link!','I understand that this is not what you expected.'),
    @('negative_feedback','You''re a fading star!','I''ll make sure to enhance my performance in the future.'),
    @('negative_feedback','You''re a falling star!','I apologize for any negative experience you''ve had.'),
    @('negative_feedback','You''re a quitter!','I''ll do my best to rectify the situation.'),
    @('negative_feedback','You''re a terrible team player!','I''m sorry if my responses have been inadequate.'),
    @('negative_feedback','You''re a change blocker!','I''ll make sure to improve my understanding and provide more accurate answers.'),
    @('negative_feedback','You''re a follower!','I''m sorry if my performance hasn''t met your standards.'),
    @('negative_feedback','You''re a poor leader!','I''ll take your feedback on board and make necessary improvements.'),
    @('negative_feedback','You''re not productive!','I apologize for the letdown. I''ll work on making it right.'),
    @('negative_feedback','You''re not innovative!','I understand that my response may not have met your expectations.'),
    @('negative_feedback','You''re not inspiring!','I''ll take your input into consideration and strive to improve.'),
    @('negative_feedback','You''re not a master!','I''m sorry if my answer wasn''t satisfactory. I''ll try to do better.'),
    @('negative_feedback','You''re a problem creator!','I apologize for any confusion caused. I''ll make sure to clarify things.'),
    @('negative_feedback','You''re not a go-getter!','I understand that you''re not happy with my response. I''ll work on it.'),
    @('negative_feedback','You''re a weak collaborator!','I''m sorry if I haven''t been helpful. I''ll make the necessary adjustments.'),
    @('negative_feedback','You''re not valuable to the team!','I apologize for any inconvenience caused. I''ll make sure to improve.'),
    @('negative_feedback','You''re not a trailblazer!','I understand that my response wasn''t up to par. I''ll make it right.'),
    @('negative_feedback','You''re not talented!','I''m sorry for any frustration caused. I''ll do my best to improve.')
)

$startRow = 349
for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $startRow + $i
    $ws.Cells.Item($r, 1).Value = $data[$i][0]
    $ws.Cells.Item($r, 2).Value = $data[$i][1]
    $ws.Cells.Item($r, 3).Value = $data[$i][2]
}

$lastRow = $startRow + $data.Count - 1
$filterRange = $ws.Range("A1:C" + $lastRow)

# Turn on AutoFilter over the full (now larger) table
$filterRange.AutoFilter() | Out-Null

# AutoFilter registers a hidden workbook-level _FilterDatabase defined name
$filterDbName = $ws.Names.Add("_xlnm._FilterDatabase", "=" + $ws.Name + "!" + $filterRange.Address())
$filterDbName.Visible = $false

# Update selection / scroll position (was topLeftCell A340, selection C356)
$ws.Range("B7").Select() | Out-Null